$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "FPA001": remove the "tolerancesToggle" column (I), rename the
# (now-vacated) "surfaceTreatment" data column header to "compliances" and
# clear its values, then set new materialGroup / materialName values for the
# rows (Plastic/PEEK GF, Rubber/Moosgummi EPDM, Rubber/SBR 40).
# ---------------------------------------------------------------------------
$wsFPA001 = $wb.Worksheets.Item("FPA001")

$wsFPA001.Range("I1").EntireColumn.Delete()

$wsFPA001.Range("F1").Value = "compliances"
$wsFPA001.Range("F2").Value = ""
$wsFPA001.Range("F3").Value = ""
$wsFPA001.Range("F4").Value = ""

$wsFPA001.Range("I2").Value = "Plastic"
$wsFPA001.Range("J2").Value = "PEEK GF (natur - beige)"
$wsFPA001.Range("I3").Value = "Rubber"
$wsFPA001.Range("J3").Value = "Moosgummi EPDM (schwarz)"
$wsFPA001.Range("I4").Value = "Rubber"
$wsFPA001.Range("J4").Value = "SBR 40 (grau)"

# ---------------------------------------------------------------------------
# Sheet "FPA002-003-005-007": remove the "tolerancesToggle" column (I) too;
# no data values change besides the shift.
# ---------------------------------------------------------------------------
$wsFPA002 = $wb.Worksheets.Item("FPA002-003-005-007")
$wsFPA002.Range("I1").EntireColumn.Delete()

# Restore this sheet's remembered selection (changed from E29 to F10).
$wsFPA002.Range("F10").Select()

# ---------------------------------------------------------------------------
# Sheet "BTMI002": it is no longer the active tab; just update its
# remembered selection to C19.
# ---------------------------------------------------------------------------
$wsBTMI002 = $wb.Worksheets.Item("BTMI002")
$wsBTMI002.Range("C19").Select()

# ---------------------------------------------------------------------------
# Finally, make "FPA001" the active sheet again with its new selection
# (F14), which also clears tabSelected on every other sheet and resets the
# workbook's activeTab back to this (first) sheet.
# ---------------------------------------------------------------------------
$wsFPA001.Range("F14").Select()
